# Atualização dos preços 13/10
# Updates the futures price table on the "futuros" sheet with new quotes
# and refreshes the last-updated date. All other changed cells in the
# workbook are formulas that depend (directly or indirectly) on this
# table, so recalculating the workbook after the edit reproduces the
# rest of the diff.

$wb = $excel.ActiveWorkbook

$futuros = $wb.Worksheets.Item("futuros")

# Last updated date (13/10/2025)
$futuros.Range("D1").Value = 45943

# New KC=F quotes (column B); column C recomputes automatically via formula
$futuros.Range("B2").Value = 382.7
$futuros.Range("B3").Value = 364.55
$futuros.Range("B4").Value = 351.3
$futuros.Range("B5").Value = 340
$futuros.Range("B6").Value = 328.6
$futuros.Range("B7").Value = 320.35
$futuros.Range("B8").Value = 309.5
$futuros.Range("B9").Value = 302.3
$futuros.Range("B10").Value = 297.2
$futuros.Range("B11").Value = 291.5
$futuros.Range("B12").Value = 285.5

$excel.CalculateFullRebuild()
$excel.Calculate()

# Scroll the frozen pane on "hedge" up by one row (A16 -> A15) without
# altering the current selection (L16:L28)
$hedge = $wb.Worksheets.Item("hedge")
$hedge.Activate()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1

# Switch the active/selected sheet from "hedge" to "Sheet2"
$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet2.Activate()
